$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (column F) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 145
$wsExhibit.Range("F5").Value = 1299
$wsExhibit.Range("F6").Value = 18172
$wsExhibit.Range("F10").Value = 6849
$wsExhibit.Range("F15").Value = 64
$wsExhibit.Range("F18").Value = 1302
$wsExhibit.Range("F19").Value = 229
$wsExhibit.Range("F21").Value = 656
$wsExhibit.Range("F23").Value = 32
$wsExhibit.Range("F25").Value = 276
$wsExhibit.Range("F32").Value = 72
$wsExhibit.Range("F33").Value = 12077
$wsExhibit.Range("F35").Value = 42
$wsExhibit.Range("F36").Value = 208
$wsExhibit.Range("F37").Value = 283
$wsExhibit.Range("F38").Value = 3920

# Sheet "全部类型" (sheet4) - mirrors the same data, offset by 2 rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 145
$wsAll.Range("F5").Value = 1299
$wsAll.Range("F6").Value = 18172
$wsAll.Range("F10").Value = 6849
$wsAll.Range("F15").Value = 64
$wsAll.Range("F18").Value = 1302
$wsAll.Range("F19").Value = 229
$wsAll.Range("F21").Value = 656
$wsAll.Range("F23").Value = 32
$wsAll.Range("F25").Value = 276
$wsAll.Range("F34").Value = 72
$wsAll.Range("F35").Value = 12077
$wsAll.Range("F37").Value = 42
$wsAll.Range("F38").Value = 208
$wsAll.Range("F39").Value = 283
$wsAll.Range("F40").Value = 3920

